$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "61.825.59"
Set-TextValue $ws "E2" "  -1.87%  "
Set-TextValue $ws "D3" "2.911.76"
Set-TextValue $ws "E3" "  -2.34%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  -0.03%  "
Set-TextValue $ws "D5" "586.69"
Set-TextValue $ws "E5" "  -1.64%  "
Set-TextValue $ws "D6" "146.47"
Set-TextValue $ws "E6" "  +0.29%  "
Set-TextValue $ws "E7" "  -0.04%  "
Set-TextValue $ws "E8" "  +0.16%  "
Set-TextValue $ws "D9" "2.910.67"
Set-TextValue $ws "E9" "  -2.30%  "
Set-TextValue $ws "E10" "  -6.19%  "
Set-TextValue $ws "E11" "  +3.66%  "
Set-TextValue $ws "E12" "  -3.79%  "
Set-TextValue $ws "E13" "  +0.64%  "
Set-TextValue $ws "E14" "  -2.65%  "
Set-TextValue $ws "E15" "  -1.74%  "
Set-TextValue $ws "D16" "3.395.19"
Set-TextValue $ws "E16" "  -2.27%  "
Set-TextValue $ws "D17" "61.866.08"
Set-TextValue $ws "E17" "  -1.54%  "
Set-TextValue $ws "D18" "6.59"
Set-TextValue $ws "E18" "  -2.69%  "
Set-TextValue $ws "D19" "2.909.69"
Set-TextValue $ws "E19" "  -2.89%  "
Set-TextValue $ws "D20" "434.89"
Set-TextValue $ws "E20" "  -2.30%  "
Set-TextValue $ws "D21" "13.41"
Set-TextValue $ws "E21" "  -1.49%  "
Set-TextValue $ws "E22" "  -2.82%  "
Set-TextValue $ws "E23" "  -3.37%  "
Set-TextValue $ws "D24" "80.88"
Set-TextValue $ws "E24" "  -1.78%  "
Set-TextValue $ws "D25" "11.90"
Set-TextValue $ws "E25" "  -2.26%  "
Set-TextValue $ws "E26" "  -8.11%  "
Set-TextValue $ws "E27" "  -0.01%  "
Set-TextValue $ws "D28" "2.07"
Set-TextValue $ws "E28" "  -4.73%  "
Set-TextValue $ws "E29" "  +20.63%  "
Set-TextValue $ws "D30" "7.21"
Set-TextValue $ws "E30" "  +0.62%  "
Set-TextValue $ws "E31" "  -2.73%  "
Set-TextValue $ws "E32" "  -1.86%  "
Set-TextValue $ws "E33" "  -0.22%  "
Set-TextValue $ws "E34" "  -0.15%  "
Set-TextValue $ws "E35" "  -3.13%  "
Set-TextValue $ws "D36" "0.975"
Set-TextValue $ws "E36" "  -2.06%  "
Set-TextValue $ws "E37" "  +3.29%  "
Set-TextValue $ws "D38" "5.51"
Set-TextValue $ws "E38" "  -3.19%  "
Set-TextValue $ws "D39" "49.11"
Set-TextValue $ws "E40" "  -2.49%  "
Set-TextValue $ws "E41" "  -3.56%  "
Set-TextValue $ws "E42" "  -2.68%  "
Set-TextValue $ws "E43" "  -5.09%  "
Set-TextValue $ws "D44" "38.83"
Set-TextValue $ws "E44" "  -0.41%  "
Set-TextValue $ws "D45" "2.697.84"
Set-TextValue $ws "E45" "  -0.56%  "
Set-TextValue $ws "D46" "134.67"
Set-TextValue $ws "E46" "  -0.22%  "
Set-TextValue $ws "E47" "  -2.47%  "
Set-TextValue $ws "D48" "344.17"
Set-TextValue $ws "E48" "  -8.37%  "
Set-TextValue $ws "E49" "  +0.04%  "
Set-TextValue $ws "E50" "  -1.76%  "
